$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '34.496.90'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +1.44%  '
$c.ClearFormats()
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.787.23'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.ClearFormats()
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '222.48'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -1.71%  '
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -1.18%  '
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c.ClearFormats()
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '32.27'
$c.ClearFormats()
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +6.60%  '
$c.ClearFormats()
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.281'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.28%  '
$c.ClearFormats()
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0686'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +2.70%  '
$c.ClearFormats()
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +0.98%  '
$c.ClearFormats()
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '2.044.84'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.ClearFormats()
$c = $ws.Range('B13')
$c.NumberFormat = '@'
$c.Value = 'Chainlink'
$c.ClearFormats()
$c = $ws.Range('C13')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c.ClearFormats()
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '10.97'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +5.06%  '
$c.ClearFormats()
$c = $ws.Range('B14')
$c.NumberFormat = '@'
$c.Value = 'WrappedEther'
$c.ClearFormats()
$c = $ws.Range('C14')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.ClearFormats()
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.781.17'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -0.33%  '
$c.ClearFormats()
$c = $ws.Range('B15')
$c.NumberFormat = '@'
$c.Value = 'Polygon'
$c.ClearFormats()
$c = $ws.Range('C15')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c.ClearFormats()
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.631'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +0.72%  '
$c.ClearFormats()
$c = $ws.Range('B16')
$c.NumberFormat = '@'
$c.Value = 'WrappedBTC'
$c.ClearFormats()
$c = $ws.Range('C16')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c.ClearFormats()
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '34.514.80'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +1.52%  '
$c.ClearFormats()
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '4.28'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +2.05%  '
$c.ClearFormats()
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '68.68'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -0.71%  '
$c.ClearFormats()
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '253.73'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +0.57%  '
$c.ClearFormats()
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.0₃0781'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +5.34%  '
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.ClearFormats()
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '10.48'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +1.22%  '
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -1.62%  '
$c.ClearFormats()
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +0.20%  '
$c.ClearFormats()
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '160.63'
$c.ClearFormats()
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +1.47%  '
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -0.88%  '
$c.ClearFormats()
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.09'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -0.72%  '
$c.ClearFormats()
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -0.15%  '
$c.ClearFormats()
$c = $ws.Range('B30')
$c.NumberFormat = '@'
$c.Value = 'Hedera'
$c.ClearFormats()
$c = $ws.Range('C30')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c.ClearFormats()
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0516'
$c.ClearFormats()
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +0.09%  '
$c.ClearFormats()
$c = $ws.Range('B31')
$c.NumberFormat = '@'
$c.Value = 'Filecoin'
$c.ClearFormats()
$c = $ws.Range('C31')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.ClearFormats()
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.75'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -2.12%  '
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c.ClearFormats()
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.56'
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -0.91%  '
$c.ClearFormats()
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.87'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.71%  '
$c.ClearFormats()
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.436.26'
$c.ClearFormats()
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -4.66%  '
$c.ClearFormats()
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.640'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.63%  '
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -1.16%  '
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +2.34%  '
$c.ClearFormats()
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '85.09'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +1.78%  '
$c.ClearFormats()
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +3.16%  '
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.45%  '
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +1.48%  '
$c.ClearFormats()
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +1.73%  '
$c.ClearFormats()
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '6.02'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +4.85%  '
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -1.15%  '
$c.ClearFormats()
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0491'
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -5.10%  '
$c.ClearFormats()
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.943.95'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c.ClearFormats()
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '12.05'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +2.19%  '
$c.ClearFormats()
$c = $ws.Range('B49')
$c.NumberFormat = '@'
$c.Value = 'Quant'
$c.ClearFormats()
$c = $ws.Range('C49')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.ClearFormats()
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '103.96'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +5.88%  '
$c.ClearFormats()
$c = $ws.Range('B50')
$c.NumberFormat = '@'
$c.Value = 'PaxDollar'
$c.ClearFormats()
$c = $ws.Range('C50')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c.ClearFormats()
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -0.18%  '
$c.ClearFormats()
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '49.96'
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -3.03%  '
$c.ClearFormats()
